$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
# B2 holds the text "2" (not a number) in the source data, so force text
# entry the same way typing `'2` in the Excel UI would (quote-prefix).
$ws.Cells.Item(2,2).Value2 = "'2"  # B2
$ws.Cells.Item(2,4).Value2 = -0.0203  # D2
$ws.Cells.Item(2,7).Value2 = -0.3341621084037792  # G2
$ws.Cells.Item(2,8).Value2 = -0.8901044256588762  # H2
$ws.Cells.Item(2,9).Value2 = -1.727613143362786  # I2
$ws.Cells.Item(2,10).Value2 = -1.727613143362786  # J2
$ws.Cells.Item(2,11).Value2 = -37.85  # K2
$ws.Cells.Item(2,12).Value2 = -1.882148184982596  # L2
$ws.Cells.Item(2,21).Value2 = 23.58  # U2
$ws.Cells.Item(2,22).Value2 = 0.01495623493593809  # V2
$ws.Cells.Item(2,23).Value2 = -1.817747047925177  # W2
$ws.Cells.Item(2,24).Value2 = 0.08879919470468652  # X2
$ws.Cells.Item(2,25).Value2 = -1.906546242629863  # Y2
$ws.Cells.Item(2,26).Value2 = 0.770601375430043  # Z2
$ws.Cells.Item(2,27).Value2 = -1.148695408007903  # AA2
$ws.Cells.Item(2,28).Value2 = 0.0765788225553283  # AB2
$ws.Cells.Item(2,29).Value2 = -1.225274230563231  # AC2
$ws.Cells.Item(2,30).Value2 = 5.28  # AD2
$ws.Cells.Item(2,31).Value2 = 2.166501565128117  # AE2
$ws.Cells.Item(2,32).Value2 = 7.446501565128116  # AF2
$ws.Cells.Item(2,33).Value2 = -16.13349843487188  # AG2
$ws.Cells.Item(2,34).Value2 = 0.004700936214795809  # AH2
$ws.Cells.Item(2,35).Value2 = 0.1137625963361175  # AI2
$ws.Cells.Item(2,36).Value2 = -0.01033889443874008  # AJ2
$ws.Cells.Item(2,37).Value2 = -0.3852637596715279  # AK2
$ws.Cells.Item(2,38).Value2 = 10.933  # AL2
$ws.Cells.Item(2,39).Value2 = 10.933  # AM2
$ws.Cells.Item(2,40).Value2 = -0.1636245312838948  # AN2
$ws.Cells.Item(2,41).Value2 = -3.213207719747554  # AO2
$ws.Cells.Item(2,42).Value2 = 0.4999689620029094  # AP2
$ws.Cells.Item(2,43).Value2 = -3.213207719747554  # AQ2

# --- Row 3 updates ---
$ws.Cells.Item(3,4).Value2 = -0.0203  # D3
$ws.Cells.Item(3,7).Value2 = -0.03892215568862272  # G3
$ws.Cells.Item(3,8).Value2 = -0.2455089820359281  # H3
$ws.Cells.Item(3,9).Value2 = -0.2660060067680014  # I3
$ws.Cells.Item(3,10).Value2 = -0.2660060067680014  # J3
$ws.Cells.Item(3,11).Value2 = -5.45  # K3
$ws.Cells.Item(3,12).Value2 = -0.3263473053892216  # L3
$ws.Cells.Item(3,21).Value2 = 2.58  # U3
$ws.Cells.Item(3,22).Value2 = 0.2283185840707964  # V3
$ws.Cells.Item(3,23).Value2 = -0.6630170316301703  # W3
$ws.Cells.Item(3,24).Value2 = 0.1015336855370778  # X3
$ws.Cells.Item(3,25).Value2 = -0.7645507171672481  # Y3
$ws.Cells.Item(3,26).Value2 = 1.745674726158512  # Z3
$ws.Cells.Item(3,27).Value2 = -0.4643599630212501  # AA3
$ws.Cells.Item(3,28).Value2 = 0.07707733340887156  # AB3
$ws.Cells.Item(3,29).Value2 = -0.5414372964301216  # AC3
$ws.Cells.Item(3,30).Value2 = 3.45  # AD3
$ws.Cells.Item(3,31).Value2 = 2.166501565128117  # AE3
$ws.Cells.Item(3,32).Value2 = 5.616501565128116  # AF3
$ws.Cells.Item(3,33).Value2 = 3.036501565128116  # AG3
$ws.Cells.Item(3,34).Value2 = 0.3320131850846774  # AH3
$ws.Cells.Item(3,35).Value2 = 0.5140256038632857  # AI3
$ws.Cells.Item(3,36).Value2 = 0.2118021297827676  # AJ3
$ws.Cells.Item(3,37).Value2 = 0.3638053071019232  # AK3
$ws.Cells.Item(3,38).Value2 = 0.233  # AL3
$ws.Cells.Item(3,39).Value2 = 0.233  # AM3
$ws.Cells.Item(3,40).Value2 = -1.162007409902324  # AN3
$ws.Cells.Item(3,41).Value2 = -20.72961373390558  # AO3
$ws.Cells.Item(3,42).Value2 = -1.022735454741703  # AP3
$ws.Cells.Item(3,43).Value2 = -20.72961373390558  # AQ3

# --- Row 4 (new row) ---
$ws.Cells.Item(4,1).Value2 = "Israel"  # A4
$ws.Cells.Item(4,2).Value2 = "Nano Dimension Ltd. (NasdaqCM:NNDM)"  # B4
$ws.Cells.Item(4,3).Value2 = "Computers/Peripherals"  # C4
$ws.Cells.Item(4,7).Value2 = -1.780058651026393  # G4
$ws.Cells.Item(4,8).Value2 = -4.046920821114369  # H4
$ws.Cells.Item(4,9).Value2 = -8.885630498533724  # I4
$ws.Cells.Item(4,10).Value2 = -8.885630498533724  # J4
$ws.Cells.Item(4,11).Value2 = -32.4  # K4
$ws.Cells.Item(4,12).Value2 = -9.501466275659823  # L4
$ws.Cells.Item(4,13).Value2 = -0  # M4
$ws.Cells.Item(4,14).Value2 = -0  # N4
$ws.Cells.Item(4,15).Value2 = 0  # O4
$ws.Cells.Item(4,16).Value2 = -0  # P4
$ws.Cells.Item(4,17).Value2 = -0  # Q4
$ws.Cells.Item(4,18).Value2 = 0  # R4
$ws.Cells.Item(4,19).Value2 = 0  # S4
$ws.Cells.Item(4,21).Value2 = 21  # U4
$ws.Cells.Item(4,22).Value2 = 0.01341595860218489  # V4
$ws.Cells.Item(4,23).Value2 = -2.972477064220183  # W4
$ws.Cells.Item(4,24).Value2 = 0.07606470387229526  # X4
$ws.Cells.Item(4,25).Value2 = -3.048541768092479  # Y4
$ws.Cells.Item(4,26).Value2 = 0.206291591046582  # Z4
$ws.Cells.Item(4,27).Value2 = -1.833030852994555  # AA4
$ws.Cells.Item(4,28).Value2 = 0.07608031170178506  # AB4
$ws.Cells.Item(4,29).Value2 = -1.90911116469634  # AC4
$ws.Cells.Item(4,30).Value2 = 1.83  # AD4
$ws.Cells.Item(4,31).Value2 = 0  # AE4
$ws.Cells.Item(4,32).Value2 = 1.83  # AF4
$ws.Cells.Item(4,33).Value2 = -19.17  # AG4
$ws.Cells.Item(4,34).Value2 = 0.001167739753562244  # AH4
$ws.Cells.Item(4,35).Value2 = 0.0335595085274161  # AI4
$ws.Cells.Item(4,36).Value2 = -0.01239869868639765  # AJ4
$ws.Cells.Item(4,37).Value2 = -0.5717268118103191  # AK4
$ws.Cells.Item(4,38).Value2 = 10.7  # AL4
$ws.Cells.Item(4,39).Value2 = 10.7  # AM4
$ws.Cells.Item(4,40).Value2 = -0.06245733788395905  # AN4
$ws.Cells.Item(4,41).Value2 = -2.83177570093458  # AO4
$ws.Cells.Item(4,42).Value2 = 0.6542662116040956  # AP4
$ws.Cells.Item(4,43).Value2 = -2.83177570093458  # AQ4
